$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", matching the formatting of the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1:G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New "Save" data column values for rows 2-8
$saveValues = @(1, 0, 1, 1, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
